$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'compression men nike',
    'compression men waist',
    'compression mens',
    'compression mens 3 4',
    'compression mens running pants',
    'compression mens tights',
    'compression neck support',
    'compression nike',
    'compression pack',
    'compression pant',
    'compression pant women',
    'compression panties',
    'compression panties for women',
    'compression pants',
    'compression pants 2xu',
    'compression pants 3 4',
    'compression pants 3 pack',
    'compression pants and tops for men',
    'compression pants black',
    'compression pants boys',
    'compression pants boys basketball',
    'compression pants capri men',
    'compression pants cold weather',
    'compression pants cwx',
    'compression pants football',
    'compression pants for girls',
    'compression pants for recovery',
    'compression pants for women',
    'compression pants for women after surgery',
    'compression pants high waist women',
    'compression pants knee length',
    'compression pants long',
    'compression pants men 3/4',
    'compression pants men 3x',
    'compression pants men cwx',
    'compression pants nike',
    'compression pants plus size',
    'compression pants running',
    'compression pants size',
    'compression pants white',
    'compression pants white women',
    'compression pants with pocket',
    'compression pants women',
    'compression pants women high waist',
    'compression pants women plus size',
    'compression pants womens',
    'compression panty hose for women',
    'compression panty hose women',
    'compression pantyhose for women',
    'compression pantyhose for women for travel',
    'compression pantyhose futuro',
    'compression pantyhose women black',
    'compression pantys',
    'compression plus size',
    'compression plus size women',
    'compression pro',
    'compression recovery',
    'compression recovery pants',
    'compression recovery pants women',
    'compression recovery tights women',
    'compression ring',
    'compression running',
    'compression running capris',
    'compression running capris women',
    'compression running leg',
    'compression running leggings',
    'compression running leggings women',
    'compression running pants',
    'compression running pants men nike',
    'compression running pants with pockets',
    'compression running pants women',
    'compression running shorts',
    'compression running shorts for women',
    'compression running shorts women',
    'compression running tights',
    'compression running tights men',
    'compression running tights women',
    'compression running underwear',
    'compression running women',
    'compression sheer knee highs for women large',
    'compression shin',
    'compression short',
    'compression short leggings',
    'compression short pack',
    'compression shorts',
    'compression shorts 10 inch',
    'compression shorts 2xu',
    'compression shorts 3 4',
    'compression shorts 3 pack',
    'compression shorts 6 inch',
    'compression shorts 7 inch',
    'compression shorts basketball',
    'compression shorts big and tall',
    'compression shorts black',
    'compression shorts cheap',
    'compression shorts copper',
    'compression shorts drskin',
    'compression shorts female',
    'compression shorts for running',
    'compression shorts for swimming'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
